$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Paragraph layout before edit (1-based):
#   1 Nodes
#   2 Graph node
#   3 ;                (standalone bullet-less paragraph)
#   4 Time-Series node
#   5 ;                (standalone bullet-less paragraph)
#   6 Edges
#   7 Graph edge (to a graph node);
#   8 Virtual edge (to a Time-Series node)
#
# Target: merge each standalone ";" paragraph into the preceding
# paragraph's text (with two spaces before the semicolon) and remove
# the now-empty standalone paragraphs. Delete from the bottom up so
# earlier paragraph indices stay valid.

$tr.Paragraphs(5, 1).Delete()
$tr.Paragraphs(3, 1).Delete()

# Update "Graph node" -> "Graph node  ;"
# Use a no-common-prefix placeholder first so the engine replaces the
# run's text in place instead of splitting off a new, rPr-less run.
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "ZZZZZZZZZZZZZZZZZZZZ"
$para2.Text = "Graph node  ;"

# Update "Time-Series node" -> "Time-Series node  ;" (now paragraph 3)
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "ZZZZZZZZZZZZZZZZZZZZ"
$para3.Text = "Time-Series node  ;"
